$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '65.482.67'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '2.647.20'
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue $ws.Range('D5') '607.48'
$ws.Range('E5').Value = '  +2.44%  '
Set-TextValue $ws.Range('D6') '156.16'
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '2.652.77'
$ws.Range('E9').Value = '  +1.36%  '
Set-TextValue $ws.Range('D10') '0.124'
$ws.Range('E10').Value = '  +8.30%  '
Set-TextValue $ws.Range('D11') '5.98'
$ws.Range('E11').Value = '  +3.46%  '
Set-TextValue $ws.Range('D12') '0.404'
$ws.Range('E12').Value = '  +1.91%  '
$ws.Range('E13').Value = '  +1.48%  '
Set-TextValue $ws.Range('D14') '29.96'
$ws.Range('E14').Value = '  +5.29%  '
Set-TextValue $ws.Range('D15') '0.0000205'
$ws.Range('E15').Value = '  +19.52%  '
$ws.Range('D16').Value = '3.137.70'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('D17').Value = '65.376.92'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '2.656.85'
$ws.Range('E18').Value = '  +1.95%  '
Set-TextValue $ws.Range('D19') '12.67'
$ws.Range('E19').Value = '  +3.48%  '
Set-TextValue $ws.Range('D20') '4.89'
$ws.Range('E20').Value = '  +2.53%  '
Set-TextValue $ws.Range('D21') '358.79'
$ws.Range('E21').Value = '  +2.29%  '
Set-TextValue $ws.Range('D22') '7.45'
$ws.Range('E22').Value = '  +4.42%  '
$ws.Range('E23').Value = '  +0.04%  '
Set-TextValue $ws.Range('D24') '70.20'
$ws.Range('E24').Value = '  +3.98%  '
$ws.Range('E25').Value = '  +0.29%  '
Set-TextValue $ws.Range('D26') '9.54'
$ws.Range('E26').Value = '  +2.82%  '
Set-TextValue $ws.Range('D27') '0.0000105'
$ws.Range('E27').Value = '  +16.13%  '
Set-TextValue $ws.Range('D28') '1.63'
$ws.Range('E28').Value = '  -0.65%  '
Set-TextValue $ws.Range('D29') '0.169'
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('E30').Value = '  +8.06%  '
Set-TextValue $ws.Range('D31') '8.10'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('E32').Value = '  +0.30%  '
Set-TextValue $ws.Range('D33') '532.19'
$ws.Range('E33').Value = '  -2.12%  '
Set-TextValue $ws.Range('D34') '1.78'
$ws.Range('E34').Value = '  -1.26%  '
Set-TextValue $ws.Range('D35') '5.51'
$ws.Range('E35').Value = '  -1.85%  '
Set-TextValue $ws.Range('D36') '6.38'
$ws.Range('E36').Value = '  +2.97%  '
Set-TextValue $ws.Range('D37') '0.432'
$ws.Range('E37').Value = '  +2.20%  '
Set-TextValue $ws.Range('D38') '20.64'
$ws.Range('E38').Value = '  +2.50%  '
Set-TextValue $ws.Range('D39') '163.11'
$ws.Range('E39').Value = '  -0.37%  '
Set-TextValue $ws.Range('D40') '1.99'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D43') '167.76'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D44') '41.95'
$ws.Range('E44').Value = '  +0.96%  '
Set-TextValue $ws.Range('D45') '4.15'
$ws.Range('E45').Value = '  +1.33%  '
Set-TextValue $ws.Range('D46') '2.32'
$ws.Range('E46').Value = '  +4.81%  '
Set-TextValue $ws.Range('D47') '0.0612'
$ws.Range('E47').Value = '  +2.78%  '
Set-TextValue $ws.Range('D48') '23.05'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D49') '0.656'
$ws.Range('E49').Value = '  +2.53%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D50') '0.0264'
$ws.Range('E50').Value = '  +5.37%  '
Set-TextValue $ws.Range('D51') '0.0981'
$ws.Range('E51').Value = '  +0.04%  '
